$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header E1 text (same shared-string slot but new text)
$ws.Range("E1").Value = "Название университетов"

# Row 2 - update existing values
$ws.Range("A2").Value = "PHYSICS"
$ws.Range("B2").Value = 4.53000020980835
$ws.Range("C2").Value = 8.0
$ws.Range("D2").Value = 2.0
$ws.Range("E2").Value = " МВУ`n МПИ`n"

# Row 3 - new
$ws.Range("A3").Value = "MEDICINE"
$ws.Range("B3").Value = 4.329999923706055
$ws.Range("C3").Value = 3.0
$ws.Range("D3").Value = 3.0
$ws.Range("E3").Value = " МГМУ`n ТУМ`n СМИ`n"

# Row 4 - new
$ws.Range("A4").Value = "LINGUISTICS"
$ws.Range("B4").Value = 0.0
$ws.Range("C4").Value = 0.0
$ws.Range("D4").Value = 1.0
$ws.Range("E4").Value = " ВЛПУ`n"

# Row 5 - new
$ws.Range("A5").Value = "MATHEMATICS"
$ws.Range("B5").Value = 0.0
$ws.Range("C5").Value = 0.0
$ws.Range("D5").Value = 1.0
$ws.Range("E5").Value = " КУВ`n"
